$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.72000000000058"
$ws.Range("H2").Value = [double]"4.829470157119431e-13"
$ws.Range("I2").Value = [double]"4.829470157119431e-13"
$ws.Range("L2").Value = [double]"52.54543715425431"
$ws.Range("M2").Value = "[41.9470610247796, 63.14381328372902]"
$ws.Range("N2").Value = [double]"5.437872374614017e-13"
$ws.Range("O2").Value = [double]"5.437872374614017e-13"
$ws.Range("P2").Value = [double]"1.566079220708426"
$ws.Range("Q2").Value = "[1.3270791789938876, 1.8050792624229643]"
$ws.Range("T2").Value = [double]"49.6360032884629"
$ws.Range("U2").Value = "[42.37845452967621, 56.89355204724959]"
$ws.Range("X2").Value = [double]"19.30930930930974"
$ws.Range("Y2").Value = [double]"18.33097097097139"
$ws.Range("Z2").Value = [double]"20.2876476476481"

# Row 3
$ws.Range("F3").Value = [double]"25.72000000000058"
$ws.Range("H3").Value = [double]"2.851052727237402e-13"
$ws.Range("I3").Value = [double]"2.851052727237402e-13"
$ws.Range("L3").Value = [double]"56.51799847772477"
$ws.Range("M3").Value = "[42.70785147740317, 70.32814547804638]"
$ws.Range("N3").Value = [double]"1.532345361709986e-10"
$ws.Range("O3").Value = [double]"1.532345361709986e-10"
$ws.Range("P3").Value = [double]"1.17613178422681"
$ws.Range("Q3").Value = "[0.9245527929483472, 1.4277107755052736]"
$ws.Range("R3").Value = [double]"3.309352791802667e-12"
$ws.Range("S3").Value = [double]"3.309352791802667e-12"
$ws.Range("T3").Value = [double]"52.8119865629342"
$ws.Range("U3").Value = "[45.29637261721898, 60.327600508649425]"
$ws.Range("X3").Value = [double]"20.90554554554601"
$ws.Range("Y3").Value = [double]"19.87571571571615"
$ws.Range("Z3").Value = [double]"21.93537537537587"

# Row 4
$ws.Range("F4").Value = [double]"25.72000000000058"
$ws.Range("H4").Value = [double]"7.069678176208072e-12"
$ws.Range("I4").Value = [double]"7.069678176208072e-12"
$ws.Range("L4").Value = [double]"59.47036146263922"
$ws.Range("M4").Value = "[43.085534949204494, 75.85518797607395]"
$ws.Range("N4").Value = [double]"3.538670245717412e-09"
$ws.Range("O4").Value = [double]"3.538670245717412e-09"
$ws.Range("P4").Value = [double]"0.6981317007977319"
$ws.Range("Q4").Value = "[0.40881586082750054, 0.9874475407679633]"
$ws.Range("R4").Value = [double]"1.463451314398512e-05"
$ws.Range("S4").Value = [double]"1.463451314398512e-05"
$ws.Range("T4").Value = [double]"54.28332220433702"
$ws.Range("U4").Value = "[45.82577892107593, 62.74086548759811]"
$ws.Range("V4").Value = [double]"2.220446049250313e-16"
$ws.Range("W4").Value = [double]"2.220446049250313e-16"
$ws.Range("X4").Value = [double]"22.86222222222274"
$ws.Range("Y4").Value = [double]"21.67791791791841"
$ws.Range("Z4").Value = [double]"24.04652652652707"

# Row 5
$ws.Range("F5").Value = [double]"25.72000000000058"
$ws.Range("H5").Value = [double]"2.394751064116463e-13"
$ws.Range("I5").Value = [double]"2.394751064116463e-13"
$ws.Range("L5").Value = [double]"58.69092377310272"
$ws.Range("M5").Value = "[45.5391832753797, 71.84266427082574]"
$ws.Range("N5").Value = [double]"1.316968756270853e-11"
$ws.Range("O5").Value = [double]"1.316968756270853e-11"
$ws.Range("P5").Value = [double]"0.3207632138800394"
$ws.Range("Q5").Value = "[0.08176317216550011, 0.5597632555945786]"
$ws.Range("R5").Value = [double]"0.009656747921369391"
$ws.Range("S5").Value = [double]"0.009656747921369391"
$ws.Range("T5").Value = [double]"56.58976503104024"
$ws.Range("U5").Value = "[49.197500317642564, 63.98202974443791]"
$ws.Range("X5").Value = [double]"24.40696696696752"
$ws.Range("Y5").Value = [double]"23.42862862862916"
$ws.Range("Z5").Value = [double]"25.38530530530588"

# Row 6
$ws.Range("F6").Value = [double]"25.72000000000058"
$ws.Range("H6").Value = [double]"2.486510997101732e-11"
$ws.Range("I6").Value = [double]"2.486510997101732e-11"
$ws.Range("L6").Value = [double]"54.9106333669156"
$ws.Range("M6").Value = "[41.55523373252643, 68.26603300130478]"
$ws.Range("N6").Value = [double]"1.349176326215229e-10"
$ws.Range("O6").Value = [double]"1.349176326215229e-10"
$ws.Range("P6").Value = [double]"-0.1257894956392311"
$ws.Range("Q6").Value = "[-0.37736848691769254, 0.12578949563923025]"
$ws.Range("R6").Value = [double]"0.3192945136149032"
$ws.Range("S6").Value = [double]"0.3192945136149032"
$ws.Range("T6").Value = [double]"53.53068780074351"
$ws.Range("U6").Value = "[45.409125989735344, 61.652249611751685]"
$ws.Range("X6").Value = [double]"0.5149149149149288"
$ws.Range("Y6").Value = [double]"-0.5149149149149221"
$ws.Range("Z6").Value = [double]"1.54474474474478"

# Row 7
$ws.Range("F7").Value = [double]"25.72000000000058"
$ws.Range("H7").Value = [double]"8.526512829121202e-14"
$ws.Range("I7").Value = [double]"8.526512829121202e-14"
$ws.Range("L7").Value = [double]"54.15146002750111"
$ws.Range("M7").Value = "[41.143159910163064, 67.15976014483915]"
$ws.Range("N7").Value = [double]"9.570189085650327e-11"
$ws.Range("O7").Value = [double]"9.570189085650327e-11"
$ws.Range("P7").Value = [double]"-0.4905790329930007"
$ws.Range("Q7").Value = "[-0.7170001251436169, -0.2641579408423844]"
$ws.Range("R7").Value = [double]"7.382293980140453e-05"
$ws.Range("S7").Value = [double]"7.382293980140453e-05"
$ws.Range("T7").Value = [double]"53.42897010038467"
$ws.Range("U7").Value = "[46.565700409605746, 60.2922397911636]"
$ws.Range("X7").Value = [double]"2.008168168168215"
$ws.Range("Y7").Value = [double]"1.081321321321345"
$ws.Range("Z7").Value = [double]"2.935015015015085"

# Row 8
$ws.Range("F8").Value = [double]"25.72000000000058"
$ws.Range("H8").Value = [double]"1.110223024625157e-14"
$ws.Range("I8").Value = [double]"1.110223024625157e-14"
$ws.Range("L8").Value = [double]"57.39962132096292"
$ws.Range("M8").Value = "[44.01358064792052, 70.78566199400532]"
$ws.Range("N8").Value = [double]"4.161093691834594e-11"
$ws.Range("O8").Value = [double]"4.161093691834594e-11"
$ws.Range("P8").Value = [double]"-0.9937370155499243"
$ws.Range("Q8").Value = "[-1.2327370572644627, -0.754736973835386]"
$ws.Range("R8").Value = [double]"9.891532037897832e-11"
$ws.Range("S8").Value = [double]"9.891532037897832e-11"
$ws.Range("T8").Value = [double]"51.56772176347614"
$ws.Range("U8").Value = "[44.504738535611835, 58.63070499134044]"
$ws.Range("X8").Value = [double]"4.067827827827919"
$ws.Range("Y8").Value = [double]"3.089489489489561"
$ws.Range("Z8").Value = [double]"5.046166166166278"

# Row 9
$ws.Range("F9").Value = [double]"24.92000000000046"
$ws.Range("H9").Value = [double]"1.867250798426312e-11"
$ws.Range("I9").Value = [double]"1.867250798426312e-11"
$ws.Range("L9").Value = [double]"50.58423657308949"
$ws.Range("M9").Value = "[37.22764638409616, 63.94082676208282]"
$ws.Range("N9").Value = [double]"1.206297506328724e-09"
$ws.Range("O9").Value = [double]"1.206297506328724e-09"
$ws.Range("P9").Value = [double]"-1.333368653775848"
$ws.Range("Q9").Value = "[-1.6226844937460783, -1.0440528138056173]"
$ws.Range("R9").Value = [double]"5.081712828314267e-12"
$ws.Range("S9").Value = [double]"5.081712828314267e-12"
$ws.Range("T9").Value = [double]"53.30481448996535"
$ws.Range("U9").Value = "[45.56304540529881, 61.04658357463188]"
$ws.Range("X9").Value = [double]"5.288328328328422"
$ws.Range("Y9").Value = [double]"4.140860860860936"
$ws.Range("Z9").Value = [double]"6.435795795795908"

# Row 10
$ws.Range("F10").Value = [double]"24.92000000000046"
$ws.Range("H10").Value = [double]"3.885780586188048e-15"
$ws.Range("I10").Value = [double]"3.885780586188048e-15"
$ws.Range("L10").Value = [double]"59.05643115880086"
$ws.Range("M10").Value = "[48.39271904303304, 69.72014327456868]"
$ws.Range("N10").Value = [double]"1.509903313490213e-14"
$ws.Range("O10").Value = [double]"1.509903313490213e-14"
$ws.Range("P10").Value = [double]"-1.471737098979002"
$ws.Range("Q10").Value = "[-1.6730002920017712, -1.2704739059562327]"
$ws.Range("R10").Value = [double]"0"
$ws.Range("S10").Value = [double]"0"
$ws.Range("T10").Value = [double]"52.19867591157025"
$ws.Range("U10").Value = "[45.24004930515436, 59.15730251798614]"
$ws.Range("X10").Value = [double]"5.837117117117224"
$ws.Range("Y10").Value = [double]"5.038878878878972"
$ws.Range("Z10").Value = [double]"6.635355355355475"

# Row 11
$ws.Range("F11").Value = [double]"24.92000000000046"
$ws.Range("H11").Value = [double]"1.110223024625157e-16"
$ws.Range("I11").Value = [double]"1.110223024625157e-16"
$ws.Range("L11").Value = [double]"59.97471461119393"
$ws.Range("M11").Value = "[50.627895784562895, 69.32153343782497]"
$ws.Range("N11").Value = [double]"2.220446049250313e-16"
$ws.Range("O11").Value = [double]"2.220446049250313e-16"
$ws.Range("P11").Value = [double]"-1.647842392873925"
$ws.Range("Q11").Value = "[-1.8239476867688484, -1.471737098979002]"
$ws.Range("R11").Value = [double]"0"
$ws.Range("S11").Value = [double]"0"
$ws.Range("T11").Value = [double]"56.16150624903397"
$ws.Range("U11").Value = "[50.05070573613799, 62.27230676192995]"
$ws.Range("X11").Value = [double]"6.535575575575695"
$ws.Range("Y11").Value = [double]"5.837117117117224"
$ws.Range("Z11").Value = [double]"7.234034034034166"
